$wb = $excel.ActiveWorkbook

# --- Sheets ---
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- New row 3 on "Logs" ---
$answerText = @"
Beste afzender,
Bedankt voor je interesse in de VentiQ-250. Helaas kan ik de datasheet niet direct via e-mail versturen, maar je kunt de datasheet vinden op onze website [link naar datasheet]. Mocht je nog vragen hebben of meer informatie nodig hebben, dan hoor ik het graag!
Met vriendelijke groet,
[Naam] 
E-mailassistent van [Bedrijfsnaam]
"@

$logs.Range("A3").Value = "Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #13: Kun je mij de datasheet van de VentiQ-250 sturen?"
$logs.Range("D3").Value = "Productinformatie"
$logs.Range("E3").Value = $answerText
$logs.Range("F3").Value = "2025-08-02 00:02:59"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Nee"
$logs.Range("I3").Value = "Ja"
$logs.Range("J3").Value = "Nee"

# Writing the multi-line answer text triggers an auto row-height bump;
# AutoFit brings row 3 back to the sheet's default (unset) row height,
# matching row 2's untouched <row> element (no ht/customHeight attrs).
$logs.Rows.Item(3).AutoFit()

# --- Extend conditional formatting ranges down to row 3 ---
$logs.Range("D2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D3"))
$logs.Range("G2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G3"))
$logs.Range("H2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H3"))
$logs.Range("I2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I3"))
$logs.Range("J2").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J3"))

# --- New row 3 on "Dashboard" ---
$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 1

# --- Update chart series to include the new Dashboard row ---
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
